$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$ws.Range("D2").Value = "29.104.88"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.835.08"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.31%  "
Set-TextValue "D5" "243.83"
$ws.Range("E5").Value = "  +0.66%  "
Set-TextValue "D6" "0.6278"
$ws.Range("E6").Value = "  +0.23%  "
Set-TextValue "D7" "1.002"
$ws.Range("E7").Value = "  +0.26%  "
Set-TextValue "D8" "0.07518"
$ws.Range("E8").Value = "  -1.21%  "
Set-TextValue "D9" "0.2925"
$ws.Range("E9").Value = "  -0.06%  "
Set-TextValue "D10" "23.20"
$ws.Range("E10").Value = "  +2.84%  "
Set-TextValue "D11" "0.07696"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "1.835.09"
$ws.Range("E12").Value = "  -0.03%  "
Set-TextValue "D13" "5.004"
$ws.Range("E13").Value = "  +1.08%  "
Set-TextValue "D14" "0.6683"
$ws.Range("E14").Value = "  +0.37%  "
Set-TextValue "D15" "82.67"
$ws.Range("E15").Value = "  -0.02%  "
Set-TextValue "D16" "0.000009378"
$ws.Range("E16").Value = "  -8.34%  "
Set-TextValue "D17" "5.997"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "29.116.16"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "2.078.32"
$ws.Range("E19").Value = "  -0.28%  "
Set-TextValue "D20" "12.59"
$ws.Range("E20").Value = "  +1.99%  "
Set-TextValue "D21" "223.60"
$ws.Range("E21").Value = "  -1.23%  "
Set-TextValue "D22" "1.005"
$ws.Range("E22").Value = "  +0.67%  "
Set-TextValue "D23" "7.108"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("E24").Value = "  +0.33%  "
Set-TextValue "D25" "160.04"
$ws.Range("E25").Value = "  +1.17%  "
Set-TextValue "D26" "0.1394"
$ws.Range("E26").Value = "  +1.40%  "
Set-TextValue "D27" "8.504"
$ws.Range("E27").Value = "  +0.23%  "
Set-TextValue "D28" "17.92"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E29").Value = "  +0.86%  "
Set-TextValue "D30" "0.05691"
$ws.Range("E30").Value = "  +9.13%  "
$ws.Range("E31").Value = "  +1.28%  "
Set-TextValue "D32" "4.069"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("E33").Value = "  +1.36%  "
Set-TextValue "D34" "0.7450"
$ws.Range("E34").Value = "  +1.21%  "
Set-TextValue "D35" "1.844"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +0.07%  "
Set-TextValue "D37" "2.671"
$ws.Range("E37").Value = "  -1.11%  "
Set-TextValue "D38" "2.764"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").Value = "1.220.02"
$ws.Range("E39").Value = "  -1.75%  "
Set-TextValue "D40" "0.01783"
$ws.Range("E40").Value = "  -0.08%  "
Set-TextValue "D41" "6.526"
$ws.Range("E41").Value = "  +2.92%  "
Set-TextValue "D42" "0.8924"
$ws.Range("E42").Value = "  -0.20%  "
Set-TextValue "D43" "1.002"
$ws.Range("E43").Value = "  +0.28%  "
Set-TextValue "D44" "101.96"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.976.84"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D46" "0.00000000125"
$ws.Range("E46").Value = "  +0.92%  "
Set-TextValue "D47" "65.73"
$ws.Range("E47").Value = "  +2.28%  "
Set-TextValue "D48" "0.07657"
$ws.Range("E48").Value = "  +11.61%  "
Set-TextValue "D49" "0.5094"
$ws.Range("E49").Value = "  -0.18%  "
Set-TextValue "D50" "0.4081"
$ws.Range("E50").Value = "  +1.18%  "
Set-TextValue "D51" "9.013"
$ws.Range("E51").Value = "  +1.84%  "
